$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Remove the "Meta description" paragraph that currently follows the
#    H1 title ("Play Alexander's Conquest Free - Online Slot Game Review").
# ------------------------------------------------------------------
$metaPara = $d.Paragraphs.Item(2)
if ($metaPara.Range.Text -like "Meta description*") {
    $metaPara.Range.Delete()
}

# ------------------------------------------------------------------
# 2) At the very end of the document, the old "Prompt: ..." paragraph is
#    replaced by two paragraphs:
#      - a bold paragraph repeating the title
#      - an italic paragraph with the former meta-description text
# ------------------------------------------------------------------
$count = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($count)

$titleText = "Play Alexander’s Conquest Free - Online Slot Game Review"
$descText = "Experience the epic online slot game Alexander’s Conquest. Play for free and read our review on the shifting game grid, graphics, bonuses, RTP, and more."

$xml = "<w:p xmlns:w=`"http://schemas.openxmlformats.org/wordprocessingml/2006/main`"><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>$titleText</w:t></w:r></w:p><w:p xmlns:w=`"http://schemas.openxmlformats.org/wordprocessingml/2006/main`"><w:r/><w:r><w:rPr><w:i/></w:rPr><w:t>$descText</w:t></w:r></w:p>"

$lastPara.Range.InsertXML($xml)
